# Applies the "Implemented sync case server side" edit:
#  - Para 2 ("Hash del file...")      -> split into proofErr-wrapped runs
#  - Para 4 ("Come fare la sync?")    -> split "fare la sync?" + proofErr on "sync"
#  - Para 5 ("Se lato server...")     -> proofErr-wrap "path-hash" and "json"
#  - Para 7 (was "Nel caso in cui...") -> replaced with new text about
#    do_write/enqueue_msg, underline formatting removed, proofErr-wrapped
#    jargon words.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p2xml = '<w:p ' + $wNs + '>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Hash</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> del file: </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>path</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>, dimensione, ultima modifica</w:t></w:r>' + `
    '</w:p>'

$p4xml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:t xml:space="preserve">Come </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">fare la </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>sync</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>?</w:t></w:r>' + `
    '</w:p>'

$p5xml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:t xml:space="preserve">Se lato server ci salvassimo la struttura dati contente </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>path-hash</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> in un </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>json</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> salvato in corrispondenza di username-password, ogni volta che il client si riconnette possiamo riprendere l’intera mappa e, nello switch, confrontarla con la mappa lato client per verificare se e quali file sono stati modificati.</w:t></w:r>' + `
    '</w:p>'

$p7xml = '<w:p ' + $wNs + '>' + `
    '<w:r><w:t xml:space="preserve">Bisogna gestire la scrittura, perché la </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>do_write</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> è specifica di ogni </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>server_session</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> mentre la </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>enqueue_msg</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> è </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>comune.</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>'

$d.Paragraphs(2).Range.InsertXML($p2xml)
$d.Paragraphs(4).Range.InsertXML($p4xml)
$d.Paragraphs(5).Range.InsertXML($p5xml)
$d.Paragraphs(7).Range.InsertXML($p7xml)
